# Don't export id field in CSV/Export
# Remove the two rows that hold the "id" parameter (name/value pairs) from
# the "example" worksheet. Deleting the entire row shifts everything below
# up, which is what the target workbook shows (rows 17 and 31 in the
# original numbering held "id" / 1 and "id" / 2 respectively).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First "id" row (originally row 17: A17="id", B17=1)
$ws.Rows(17).Select()
$ws.Rows(17).Delete()

# Second "id" row - after the first deletion the row that used to be 31
# is now row 30 (originally row 31: A31="id", B31=2)
$ws.Rows(30).Select()
$ws.Rows(30).Delete()
